$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set G2 (Background Processing) to TRUE
$ws.Range("G2").Value = $true

# Delete row 3 (duplicate row, now redundant since G2 is TRUE)
$ws.Rows(3).Delete()

# Update selection to match target state (row 2 selected)
$ws.Range("A2:XFD2").Select()
